$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.019597685490375
$ws.Range("D2").Value = 1.024733527851257
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.030629225673865
$ws.Range("I2").Value = 1.028830943338726
$ws.Range("J2").Value = 1.02479954828536
$ws.Range("K2").Value = 1.02756116299202
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.033439683782387
$ws.Range("N2").Value = 1.012300903115731

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.020714901171515
$ws.Range("D3").Value = 1.025527507425509
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.031925499853288
$ws.Range("I3").Value = 1.029029899356856
$ws.Range("J3").Value = 1.025552711462865
$ws.Range("K3").Value = 1.028162332178411
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.034543043100339
$ws.Range("N3").Value = 1.012554411321746

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.021437401891201
$ws.Range("D4").Value = 1.026040778061286
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.032764067111685
$ws.Range("I4").Value = 1.029157079568751
$ws.Range("J4").Value = 1.026039144722261
$ws.Range("K4").Value = 1.028550192291059
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.035256242790788
$ws.Range("N4").Value = 1.012718022256808

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021741044002686
$ws.Range("D5").Value = 1.026256440608719
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.033116553126514
$ws.Range("I5").Value = 1.029210173364455
$ws.Range("J5").Value = 1.026243423322011
$ws.Range("K5").Value = 1.02871297700514
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.035555894903301
$ws.Range("N5").Value = 1.012786702538153

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021792021195527
$ws.Range("D6").Value = 1.026292644478228
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.033175734352224
$ws.Range("I6").Value = 1.02921906619114
$ws.Range("J6").Value = 1.026277709840511
$ws.Range("K6").Value = 1.028740293355596
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.035606197486032
$ws.Range("N6").Value = 1.012798228298006

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.021441459551891
$ws.Range("D7").Value = 1.026043660212024
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.032768777231438
$ws.Range("I7").Value = 1.029157790475195
$ws.Range("J7").Value = 1.026041875158292
$ws.Range("K7").Value = 1.028552368495376
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.035260247450217
$ws.Range("N7").Value = 1.012718940365688

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.019975340124612
$ws.Range("D8").Value = 1.025001957986375
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.031067352242587
$ws.Range("I8").Value = 1.028898504057947
$ws.Range("J8").Value = 1.025054272970672
$ws.Range("K8").Value = 1.027764566515293
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.033812725183149
$ws.Range("N8").Value = 1.012386665684877

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.017388620050774
$ws.Range("D9").Value = 1.023162603228984
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.028067519987127
$ws.Range("I9").Value = 1.028429672200316
$ws.Range("J9").Value = 1.023306957869803
$ws.Range("K9").Value = 1.026367634904891
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.031256193543481
$ws.Range("N9").Value = 1.011797885236232

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.015661866403524
$ws.Range("D10").Value = 1.021933832517461
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.026066337092497
$ws.Range("I10").Value = 1.028109082132823
$ws.Range("J10").Value = 1.022137298866519
$ws.Range("K10").Value = 1.025430452318782
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.029547807009821
$ws.Range("N10").Value = 1.011403153188887

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.014913600269527
$ws.Range("D11").Value = 1.021401155151256
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.025199461761368
$ws.Range("I11").Value = 1.027968354435401
$ws.Range("J11").Value = 1.021629676711034
$ws.Range("K11").Value = 1.025023236407134
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.028807072006167
$ws.Range("N11").Value = 1.01123170200271

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.014635573085454
$ws.Range("D12").Value = 1.021203202641682
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.024877410065687
$ws.Range("I12").Value = 1.027915794632692
$ws.Range("J12").Value = 1.021440948967856
$ws.Range("K12").Value = 1.024871765770852
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.028531778401536
$ws.Range("N12").Value = 1.011167937511177

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.014695214868767
$ws.Range("D13").Value = 1.021245668317873
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.024946493812384
$ws.Range("I13").Value = 1.02792708189324
$ws.Range("J13").Value = 1.021481439614735
$ws.Range("K13").Value = 1.024904266365068
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.028590836717492
$ws.Range("N13").Value = 1.0111816188313

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.014890620259707
$ws.Range("D14").Value = 1.021384794223534
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.025172842019509
$ws.Range("I14").Value = 1.027964015687688
$ws.Range("J14").Value = 1.021614079978899
$ws.Range("K14").Value = 1.025010720136579
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.028784319259199
$ws.Range("N14").Value = 1.011226432841484

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.015011004268752
$ws.Range("D15").Value = 1.021470502023064
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.02531229516207
$ws.Range("I15").Value = 1.027986733742282
$ws.Range("J15").Value = 1.021695780956546
$ws.Range("K15").Value = 1.025076281636341
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.028903510068185
$ws.Range("N15").Value = 1.011254033636759

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.015711514163188
$ws.Range("D16").Value = 1.021969171629501
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.026123861122484
$ws.Range("I16").Value = 1.0281183814977
$ws.Range("J16").Value = 1.022170963703067
$ws.Range("K16").Value = 1.025457448135599
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.029596946026962
$ws.Range("N16").Value = 1.011414520654054

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.016150771299191
$ws.Range("D17").Value = 1.022281809793674
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.026632839189161
$ws.Range("I17").Value = 1.028200449061402
$ws.Range("J17").Value = 1.022468724005818
$ws.Range("K17").Value = 1.025696165890567
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.030031652602784
$ws.Range("N17").Value = 1.01151504790507

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.016406927661938
$ws.Range("D18").Value = 1.022464107377607
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.026929683924532
$ws.Range("I18").Value = 1.028248133407969
$ws.Range("J18").Value = 1.022642291472793
$ws.Range("K18").Value = 1.025835269955114
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.030285114017121
$ws.Range("N18").Value = 1.011573632672279

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.016494261115729
$ws.Range("D19").Value = 1.022526256157089
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.027030894707727
$ws.Range("I19").Value = 1.02826436129079
$ws.Range("J19").Value = 1.022701454688469
$ws.Range("K19").Value = 1.025882677806699
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.03037152167539
$ws.Range("N19").Value = 1.011593599909976

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.016103648858498
$ws.Range("D20").Value = 1.022248272780483
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.026578234152237
$ws.Range("I20").Value = 1.028191663051007
$ws.Range("J20").Value = 1.022436788641168
$ws.Range("K20").Value = 1.025670567803059
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.029985022605253
$ws.Range("N20").Value = 1.011504267567844

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.014833080683969
$ws.Range("D21").Value = 1.021343827667533
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.025106189673291
$ws.Range("I21").Value = 1.02795314753827
$ws.Range("J21").Value = 1.021575025517961
$ws.Range("K21").Value = 1.024979378026516
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.02872734767973
$ws.Range("N21").Value = 1.011213238440307

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.014033714853706
$ws.Range("D22").Value = 1.020774632148781
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.024180333692351
$ws.Range("I22").Value = 1.027801521161658
$ws.Range("J22").Value = 1.021032191958409
$ws.Range("K22").Value = 1.024543569665774
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.027935718906129
$ws.Range("N22").Value = 1.011029794513764

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.014457522783219
$ws.Range("D23").Value = 1.021076424356094
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.024671179114496
$ws.Range("I23").Value = 1.027882058832386
$ws.Range("J23").Value = 1.021320054327127
$ws.Range("K23").Value = 1.024774716750804
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.028355460436384
$ws.Range("N23").Value = 1.011127085544578

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.016124941630124
$ws.Range("D24").Value = 1.022263426895115
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.026602907920895
$ws.Range("I24").Value = 1.028195633640167
$ws.Range("J24").Value = 1.022451219200665
$ws.Range("K24").Value = 1.025682134896784
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.030006092988612
$ws.Range("N24").Value = 1.011509138895914

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.018057742446715
$ws.Range("D25").Value = 1.023638565833172
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.028843265777851
$ws.Range("I25").Value = 1.028552292012887
$ws.Range("J25").Value = 1.023759519693215
$ws.Range("K25").Value = 1.026729812033044
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.031917819807296
$ws.Range("N25").Value = 1.011950487952215
